# calc_load_and_loc.xlsx — span 2 test data completion
#
# Moves the "cur axle location" (C8:C25) values 50 further down-span and
# updates the span boundaries (x / begin span / end span in C2:C4) so the
# calc_reactions / calc_pier_reactions tests exercise span 2 (x=125).
# Dependent formulas in D:I, row 27, row 29 and B31:C36 recalculate
# automatically. Also moves the active selection to D22 and collapses the
# now-identical "xt/xl/xr" number format onto the "Pt/Pl/Pr" one
# (0.000 -> 0.000000) for C31:C36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- span boundaries -------------------------------------------------
$ws.Range("C2").Value = 125   # x
$ws.Range("C3").Value = 100   # begin span
$ws.Range("C4").Value = 150   # end span

# --- axle locations (each +50, shifting the load train into span 2) --
$ws.Range("C8").Value  = 181
$ws.Range("C9").Value  = 173
$ws.Range("C10").Value = 168
$ws.Range("C11").Value = 163
$ws.Range("C12").Value = 158
$ws.Range("C13").Value = 149
$ws.Range("C14").Value = 144
$ws.Range("C15").Value = 138
$ws.Range("C16").Value = 133
$ws.Range("C17").Value = 125
$ws.Range("C18").Value = 117
$ws.Range("C19").Value = 112
$ws.Range("C20").Value = 107
$ws.Range("C21").Value = 102
$ws.Range("C22").Value = 93
$ws.Range("C23").Value = 88
$ws.Range("C24").Value = 82
$ws.Range("C25").Value = 77

# --- number format cleanup -------------------------------------------
# C32:C36 ("xt","xl","xr" etc.) used a distinct "0.0000" format from
# C31 ("Pt","0.000"). First fold them onto the same format as C31 so
# they share one style, then widen that shared format to "0.000000".
$ws.Range("C32:C36").NumberFormat = "0.000"
$ws.Range("C31:C36").NumberFormat = "0.000000"

# --- selection ----------------------------------------------------------
$ws.Range("D22").Select()
